# Crypto tracker refresh per commit "Updated cryptos list on Fri Aug  4 11:30:01 UTC 2023 with GitHub Actions".
# The sheet is scraped crypto-ranking data (Coin / Link / Price / Volume(1h)); this run just
# refreshes every row (2-51) with the latest scraped Coin/Link/Price/Volume(1h) text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as literal text in the source feed (e.g. "241.48"), including
# ones that look like plain decimals. Assigning such a string straight to .Value lets Excel
# auto-convert it to a binary number (losing the exact text / introducing float noise), so for
# those cells we briefly force Text format, assign the literal string, then clear the format
# again so the cell's style is unaffected and only its stored value changes.
$textProtectRefs = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textProtectRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "29.140.79"
$ws.Cells.Item(2, 5).Value = "  +0.03%  "
$ws.Cells.Item(3, 4).Value = "1.833.49"
$ws.Cells.Item(3, 5).Value = "  -0.09%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "241.48"
$ws.Cells.Item(5, 5).Value = "  +0.56%  "
$ws.Cells.Item(6, 4).Value = "0.6574"
$ws.Cells.Item(6, 5).Value = "  -1.08%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 2).Value = "Dogecoin"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(8, 4).Value = "0.07412"
$ws.Cells.Item(8, 5).Value = "  +0.85%  "
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 4).Value = "0.2924"
$ws.Cells.Item(9, 5).Value = "  -0.99%  "
$ws.Cells.Item(10, 2).Value = "Solana"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(10, 4).Value = "22.87"
$ws.Cells.Item(10, 5).Value = "  +0.63%  "
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(11, 4).Value = "0.07741"
$ws.Cells.Item(11, 5).Value = "  +0.78%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.836.05"
$ws.Cells.Item(12, 5).Value = "  +0.00%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "4.989"
$ws.Cells.Item(13, 5).Value = "  -0.60%  "
$ws.Cells.Item(14, 2).Value = "Polygon"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(14, 4).Value = "0.6671"
$ws.Cells.Item(14, 5).Value = "  -1.14%  "
$ws.Cells.Item(15, 2).Value = "Litecoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(15, 4).Value = "83.07"
$ws.Cells.Item(15, 5).Value = "  -3.54%  "
$ws.Cells.Item(16, 2).Value = "Uniswap"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(16, 4).Value = "6.108"
$ws.Cells.Item(16, 5).Value = "  -1.55%  "
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).Value = "0.000008607"
$ws.Cells.Item(17, 5).Value = "  +4.54%  "
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "29.140.76"
$ws.Cells.Item(18, 5).Value = "  +0.08%  "
$ws.Cells.Item(19, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(19, 4).Value = "2.084.66"
$ws.Cells.Item(19, 5).Value = "  -0.05%  "
$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).Value = "226.53"
$ws.Cells.Item(20, 5).Value = "  -0.94%  "
$ws.Cells.Item(21, 2).Value = "Avalanche"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(21, 4).Value = "12.45"
$ws.Cells.Item(21, 5).Value = "  -0.46%  "
$ws.Cells.Item(22, 2).Value = "Dai"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(22, 4).Value = "1.002"
$ws.Cells.Item(22, 5).Value = "  +0.19%  "
$ws.Cells.Item(23, 2).Value = "Chainlink"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(23, 4).Value = "7.107"
$ws.Cells.Item(23, 5).Value = "  -2.66%  "
$ws.Cells.Item(24, 2).Value = "BinanceUSD"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(24, 4).Value = "1.000"
$ws.Cells.Item(24, 5).Value = "  +0.00%  "
$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).Value = "161.42"
$ws.Cells.Item(25, 5).Value = "  +0.21%  "
$ws.Cells.Item(26, 2).Value = "Stellar"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(26, 4).Value = "0.1403"
$ws.Cells.Item(26, 5).Value = "  -0.96%  "
$ws.Cells.Item(27, 2).Value = "Cosmos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(27, 4).Value = "8.575"
$ws.Cells.Item(27, 5).Value = "  -1.13%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).Value = "17.97"
$ws.Cells.Item(28, 5).Value = "  -0.33%  "
$ws.Cells.Item(29, 2).Value = "PancakeSwap"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(29, 4).Value = "1.510"
$ws.Cells.Item(29, 5).Value = "  +0.51%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "4.099"
$ws.Cells.Item(30, 5).Value = "  -3.06%  "
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).Value = "4.043"
$ws.Cells.Item(31, 5).Value = "  -1.44%  "
$ws.Cells.Item(32, 2).Value = "Toncoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(32, 4).Value = "1.189"
$ws.Cells.Item(32, 5).Value = "  -0.85%  "
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).Value = "0.05273"
$ws.Cells.Item(33, 5).Value = "  -0.57%  "
$ws.Cells.Item(34, 2).Value = "LidoDAOToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(34, 4).Value = "1.865"
$ws.Cells.Item(34, 5).Value = "  +0.36%  "
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).Value = "0.7367"
$ws.Cells.Item(35, 5).Value = "  -1.35%  "
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "1.142"
$ws.Cells.Item(36, 5).Value = "  +1.14%  "
$ws.Cells.Item(37, 2).Value = "HuobiToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(37, 4).Value = "2.655"
$ws.Cells.Item(37, 5).Value = "  -0.99%  "
$ws.Cells.Item(38, 2).Value = "Maker"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(38, 4).Value = "1.304.80"
$ws.Cells.Item(38, 5).Value = "  -0.92%  "
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.01786"
$ws.Cells.Item(39, 5).Value = "  -1.12%  "
$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).Value = "2.740"
$ws.Cells.Item(40, 5).Value = "  +0.94%  "
$ws.Cells.Item(41, 2).Value = "TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(41, 4).Value = "0.9134"
$ws.Cells.Item(41, 5).Value = "  -0.97%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "6.081"
$ws.Cells.Item(42, 5).Value = "  +1.68%  "
$ws.Cells.Item(43, 2).Value = "PaxDollar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(43, 4).Value = "0.9998"
$ws.Cells.Item(43, 5).Value = "  +0.15%  "
$ws.Cells.Item(44, 2).Value = "Quant"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(44, 4).Value = "102.24"
$ws.Cells.Item(44, 5).Value = "  -1.12%  "
$ws.Cells.Item(45, 4).Value = "0.08132"
$ws.Cells.Item(45, 5).Value = "  +9.60%  "
$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "1.983.40"
$ws.Cells.Item(46, 5).Value = "  -0.06%  "
$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).Value = "0.5131"
$ws.Cells.Item(47, 5).Value = "  -0.70%  "
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).Value = "63.69"
$ws.Cells.Item(48, 5).Value = "  -0.20%  "
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).Value = "1.747"
$ws.Cells.Item(49, 5).Value = "  -0.72%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "0.05844"
$ws.Cells.Item(50, 5).Value = "  -1.51%  "
$ws.Cells.Item(51, 2).Value = "Aptos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(51, 4).Value = "6.766"
$ws.Cells.Item(51, 5).Value = "  -0.91%  "

foreach ($ref in $textProtectRefs) {
    $ws.Range($ref).ClearFormats()
}

